$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Please provide CCTV footage" paragraph: append a sentence about the
#    photograph and merge the (now superfluous) empty paragraph that used to
#    follow it into this one.
# ---------------------------------------------------------------------------
$cctvIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^\s*Please provide CCTV footage\s*$") {
        $cctvIdx = $i
        break
    }
}

if ($cctvIdx -gt 0) {
    $p = $d.Paragraphs.Item($cctvIdx)
    $r = $p.Range
    $r.Collapse(0)
    $r.InsertAfter(" and a current up to date photograph of prisoner")

    $nextP = $d.Paragraphs.Item($cctvIdx + 1)
    $nextP.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) "Please Note, only provide footage..." -> "Please note, only provide..."
#    and the closing sentence is replaced; the manual line break that used to
#    separate the two sentences is removed.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Please Note, only provide footage", $true, $false, $false, $false, $false, $true, 1, $false, "Please note, only provide footage", 2) | Out-Null

$d.Content.Find.Execute("Also, please confirm the identity of the individual and the timeframe of when they appear on the footage.", $true, $false, $false, $false, $false, $true, 1, $false, " Please also provide a current up-to-date photograph of the prisoner in order to confirm their identity.", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -match "ncident/specific timeframe requested\.") {
        $start = $p.Range.Start
        for ($off = 0; $off -lt $t.Length; $off++) {
            $code = [int][char]($t[$off])
            if ($code -eq 11) {
                $before = $t.Substring(0, $off)
                if ($before -match "ncident/specific timeframe requested\.$") {
                    $breakRange = $d.Range($start + $off, $start + $off + 1)
                    $breakRange.Delete()
                }
            }
        }
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Merge the two runs forming "...This statutory timeframe starts " +
#    "once a request ... Offender Subject Access Reques" into a single run
#    (drops the stale cached lastRenderedPageBreak sitting between them).
# ---------------------------------------------------------------------------
$mergeText = "This statutory timeframe starts once a request is received by any part of the Ministry of Justice (MoJ) including HMPPS, and not when it is received by the Offender Subject Access Reques"
$d.Content.Find.Execute($mergeText, $true, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2) | Out-Null
